$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - updated test data
$ws.Range("A2").Value = "TC_NA_001"
$ws.Range("B2").Value = "New account creation"
$ws.Range("C2").Value = "Roy"
$ws.Range("D2").Value = "Miller"
$ws.Range("E2").Value = "roy@yopmail.com"
$ws.Range("F2").Value = "roy@123"
$ws.Range("G2").Value = "roy@123"

# Row 3 - updated test data
$ws.Range("A3").Value = "TC_NA_002"
$ws.Range("B3").Value = "New account creation"
$ws.Range("C3").Value = "Toy2"
$ws.Range("D3").Value = "Miller2"
$ws.Range("E3").Value = "roy2@yopmail.com"
$ws.Range("F3").Value = "roy@123"
$ws.Range("G3").Value = "roy@123"
